# Append three new paragraphs at the end of the document, one per
# placeholder token: {{data_inizio}}, {{data_fine}}, {{num_prestazione}}.
#
# Each new paragraph is created with InsertParagraphAfter() off the
# collapsed end-of-story range, which duplicates the paragraph/run
# formatting of the document's final paragraph (itself {{num_prestazione}},
# using the Calibri/sz22/jc-both formatting already present) for the freshly
# inserted, still-empty paragraph. Setting .Text on that empty paragraph's
# range then fills in the placeholder text while keeping the inherited
# formatting intact.

$d = $word.ActiveDocument

$newParagraphValues = @("{{data_inizio}}", "{{data_fine}}", "{{num_prestazione}}")

foreach ($value in $newParagraphValues) {
    $endRange = $d.Content
    $endRange.Collapse(0)          # wdCollapseEnd
    $endRange.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newPara.Range.Text = $value
}
